# Update "想去人数" (interested count) values on the "展览" and "全部类型" sheets
# F2: 339 -> 340
# F3: 98  -> 99
# F4: 1418 -> 1422

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 340
    $ws.Range("F3").Value = 99
    $ws.Range("F4").Value = 1422
}
